# Helper: assign a value to a cell while forcing "text" storage so that
# date-looking strings (e.g. "2025-11-23") are not silently coerced into
# real dates by Excel's smart-entry heuristics. Afterwards the cell
# style is reset back to "Normal" so no stray number-format style
# lingers on the cell (keeps output equivalent to the plain inlineStr
# cells in the source file). Only used for the one column that actually
# risks auto-conversion (the plain yyyy-mm-dd "changed_day" column) -
# everything else is safe with a plain .Value assignment.
function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "snapshot" sheet: two players recovered (RETURN) and dropped out of
#    the injury list -> Волков Александр С (ЛОК) and Паник Рихард (ЛОК),
#    rows 14 and 15. Removing them shifts every following row up by two.
#    Every remaining row also gets a fresh `scraped_at` (column K)
#    timestamp from the new scrape run.
# ---------------------------------------------------------------------
$snapshot = $wb.Worksheets.Item("snapshot")

# Sanity-check before deleting, then drop the two RETURN rows as one block.
$snapshot.Range("A14:K15").EntireRow.Delete()

$newScrapedAt = @(
    "2025-11-23T01:54:07.841824+00:00",
    "2025-11-23T01:54:07.841843+00:00",
    "2025-11-23T01:54:10.356347+00:00",
    "2025-11-23T01:54:10.356365+00:00",
    "2025-11-23T01:54:13.079258+00:00",
    "2025-11-23T01:54:15.402052+00:00",
    "2025-11-23T01:54:18.145103+00:00",
    "2025-11-23T01:54:18.145133+00:00",
    "2025-11-23T01:54:18.145151+00:00",
    "2025-11-23T01:54:20.440432+00:00",
    "2025-11-23T01:54:23.106146+00:00",
    "2025-11-23T01:54:25.908887+00:00",
    "2025-11-23T01:54:28.660366+00:00",
    "2025-11-23T01:54:31.029127+00:00",
    "2025-11-23T01:54:36.170434+00:00",
    "2025-11-23T01:54:36.170465+00:00",
    "2025-11-23T01:54:38.336611+00:00",
    "2025-11-23T01:54:38.336643+00:00",
    "2025-11-23T01:54:38.336662+00:00",
    "2025-11-23T01:54:41.045630+00:00",
    "2025-11-23T01:54:41.045660+00:00",
    "2025-11-23T01:54:43.706448+00:00",
    "2025-11-23T01:54:43.706480+00:00",
    "2025-11-23T01:54:43.706502+00:00",
    "2025-11-23T01:54:43.706523+00:00",
    "2025-11-23T01:54:46.510282+00:00",
    "2025-11-23T01:54:46.510317+00:00",
    "2025-11-23T01:54:49.281171+00:00",
    "2025-11-23T01:54:49.281203+00:00",
    "2025-11-23T01:54:49.281222+00:00",
    "2025-11-23T01:54:49.281238+00:00",
    "2025-11-23T01:54:51.982443+00:00",
    "2025-11-23T01:54:51.982477+00:00",
    "2025-11-23T01:54:57.139100+00:00",
    "2025-11-23T01:54:57.139133+00:00",
    "2025-11-23T01:54:59.864667+00:00",
    "2025-11-23T01:54:59.864694+00:00"
)

for ($i = 0; $i -lt $newScrapedAt.Length; $i++) {
    $row = 2 + $i
    $snapshot.Cells.Item($row, 11).Value = $newScrapedAt[$i]
}

# ---------------------------------------------------------------------
# 2) "returned" sheet: this run's RETURN list replaces the previous one.
#    Clear the old 3 data rows and write the 2 players who came back.
# ---------------------------------------------------------------------
$returned = $wb.Worksheets.Item("returned")
$returned.Range("A2:G4").Clear()

$returnedData = @(
    @("ЛОК", "Локомотив", "Волков Александр С", "1369_ЛОК_волковалександрс", "RETURN", "2025-11-23T09:55:00.367713+08:00", "2025-11-23"),
    @("ЛОК", "Локомотив", "Паник Рихард", "1369_ЛОК_паникрихард", "RETURN", "2025-11-23T09:55:00.367713+08:00", "2025-11-23")
)

$r = 2
foreach ($dataRow in $returnedData) {
    for ($c = 0; $c -lt $dataRow.Length; $c++) {
        $colIndex = $c + 1
        if ($colIndex -eq 7) {
            # "changed_day" (G) looks like a plain date ("2025-11-23") and
            # would otherwise be auto-converted into a real date value.
            Set-TextValue $returned.Cells.Item($r, $colIndex) $dataRow[$c]
        } else {
            $returned.Cells.Item($r, $colIndex).Value = $dataRow[$c]
        }
    }
    $r++
}

# ---------------------------------------------------------------------
# 3) "new_injured" sheet: no newly-injured players this run -> clear the
#    two stale rows, leaving only the header.
# ---------------------------------------------------------------------
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Range("A2:G3").Clear()
